$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.841.06'
$ws.Range("E2").Value = '  +0.08%  '

$ws.Range("D3").Value = '2.268.61'
$ws.Range("E3").Value = '  -3.51%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '299.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.55'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.79%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.570'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.56%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.506'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -6.03%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.99'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.82%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0793'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.96%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.03'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.88%  '

$ws.Range("E13").Value = '  -1.51%  '

$ws.Range("D14").Value = '2.615.43'
$ws.Range("E14").Value = '  -3.52%  '

$ws.Range("D15").Value = '2.269.76'
$ws.Range("E15").Value = '  -3.83%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.63'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.82%  '

$ws.Range("D17").Value = '46.861.89'
$ws.Range("E17").Value = '  +0.29%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.797'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.47%  '

$ws.Range("D19").Value = '0.0₃0977'
$ws.Range("E19").Value = '  +2.20%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -9.36%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.31%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.73'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.67%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '246.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.23%  '

$ws.Range("E24").Value = '  -6.75%  '

$ws.Range("E25").Value = '  +0.17%  '

$ws.Range("E26").Value = '  -6.56%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '41.29'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.59%  '

$ws.Range("E28").Value = '  -1.67%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.54'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.05%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.03'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.33%  '

$ws.Range("E31").Value = '  +7.60%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.35'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.51%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '146.69'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.85%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.33'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.85%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0766'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.60%  '

$ws.Range("E36").Value = '  +1.41%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.115'
$ws.Range("D37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '15.66'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +11.78%  '

$ws.Range("E39").Value = '  -10.00%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.83'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.09%  '

$ws.Range("E41").Value = '  -7.30%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.07'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -11.26%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.13%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '93.60'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +15.42%  '

$ws.Range("D45").Value = '1.783.95'
$ws.Range("E45").Value = '  -1.60%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.90'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.32%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '70.87'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.94%  '

$ws.Range("E48").Value = '  -7.78%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '94.33'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.76%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.85'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.34%  '
